$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-19 down to 6-20
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value = "Maule"
$ws.Cells.Item(5, 4).Value = Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = 300000000
$ws.Cells.Item(5, 7).Value = "Espárragos"
$ws.Cells.Item(5, 8).Value = "Verde"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 2000
$ws.Cells.Item(5, 11).Value = 1800
$ws.Cells.Item(5, 12).Value = 1800
$ws.Cells.Item(5, 13).Value = 1800
$ws.Cells.Item(5, 14).Value = "`$/kilo"
$ws.Cells.Item(5, 15).Value = "Provincia de Linares"
$ws.Cells.Item(5, 16).Value = 1800
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
